$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Numeric-only updates (W column)
$ws.Range("W12").Value = 2781
$ws.Range("W13").Value = 591
$ws.Range("W14").Value = 100

$ws.Range("W18").Value = 2139
$ws.Range("W19").Value = 550
$ws.Range("W20").Value = 332

$ws.Range("W24").Value = 15128
$ws.Range("W25").Value = 3054
$ws.Range("W26").Value = 1532

# Rows 30-32: month label (V) shifts down along with updated volume (W)
$ws.Range("V30").Value = "Mar22"
$ws.Range("W30").Value = 4975

$ws.Range("V31").Value = "May22"
$ws.Range("W31").Value = 1697

$ws.Range("V32").Value = "Jul22"
$ws.Range("W32").Value = 1043

$ws.Range("W36").Value = 0
$ws.Range("W37").Value = 0
$ws.Range("W38").Value = 0

# B54/B55 volume fix (Out of Range bug fix)
$ws.Range("B54").Value = 4001
$ws.Range("B55").Value = 14
